# externalreferences.xlsx - add a #REF! style test for an external
# reference used inside a workbook-scoped defined name.
#
# The workbook already ships one external reference (rId2 ->
# externalLinks/externalLink1.xml, pointing at ExternalWorkbook.xlam) that
# is used directly from a worksheet formula. This change adds a second
# external reference (index [2] in formula notation) that is only ever
# consumed indirectly, through a new workbook-level defined name
# ("ExternalNamedRange"). A new formula on Sheet1 then sums that named
# range together with a blank local cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Workbook-scoped defined name pointing at cell C10 on the first sheet of
# the (new, second) external workbook reference.
$wb.Names.Add("ExternalNamedRange", "=[2]Sheet1!`$C`$10") | Out-Null

# The source workbook was saved with manual calculation turned on.
$wb.Application.Calculation = -4135  # xlCalculationManual

# New formula row that exercises the external named range.
$ws.Range("F15").Formula = "=SUM(ExternalNamedRange, D4)"

# Selection follows the newly entered cell.
$ws.Range("F15").Select() | Out-Null
